$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.905.99"
Set-TextValue $ws.Range("E2") "  -0.16%  "

Set-TextValue $ws.Range("D3") "1.875.45"
Set-TextValue $ws.Range("E3") "  -0.95%  "

Set-TextValue $ws.Range("D4") "0.9991"
Set-TextValue $ws.Range("E4") "  -0.16%  "

Set-TextValue $ws.Range("D5") "0.7388"
Set-TextValue $ws.Range("E5") "  -4.91%  "

Set-TextValue $ws.Range("D6") "242.40"
Set-TextValue $ws.Range("E6") "  -0.66%  "

Set-TextValue $ws.Range("D7") "0.9999"
Set-TextValue $ws.Range("E7") "  -0.04%  "

Set-TextValue $ws.Range("D8") "0.3158"
Set-TextValue $ws.Range("E8") "  +0.95%  "

Set-TextValue $ws.Range("E9") "  -0.68%  "

Set-TextValue $ws.Range("D10") "24.74"
Set-TextValue $ws.Range("E10") "  -4.32%  "

Set-TextValue $ws.Range("D11") "0.08364"

Set-TextValue $ws.Range("D12") "0.7505"
Set-TextValue $ws.Range("E12") "  -3.14%  "

Set-TextValue $ws.Range("D13") "5.425"
Set-TextValue $ws.Range("E13") "  +0.14%  "

Set-TextValue $ws.Range("D14") "1.836.10"
Set-TextValue $ws.Range("E14") "  -11.39%  "

Set-TextValue $ws.Range("D15") "92.57"
Set-TextValue $ws.Range("E15") "  -2.08%  "

Set-TextValue $ws.Range("D16") "29.928.56"
Set-TextValue $ws.Range("E16") "  -0.86%  "

Set-TextValue $ws.Range("D17") "6.077"
Set-TextValue $ws.Range("E17") "  -1.72%  "

Set-TextValue $ws.Range("B18") "BitcoinCash"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D18") "245.40"
Set-TextValue $ws.Range("E18") "  -0.30%  "

Set-TextValue $ws.Range("B19") "Avalanche"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "13.57"
Set-TextValue $ws.Range("E19") "  -2.41%  "

Set-TextValue $ws.Range("D20") "0.000007832"
Set-TextValue $ws.Range("E20") "  -0.42%  "

Set-TextValue $ws.Range("D21") "0.9993"
Set-TextValue $ws.Range("E21") "  -0.19%  "

Set-TextValue $ws.Range("D22") "2.127.29"
Set-TextValue $ws.Range("E22") "  -8.96%  "

Set-TextValue $ws.Range("D23") "8.026"
Set-TextValue $ws.Range("E23") "  -1.36%  "

Set-TextValue $ws.Range("D24") "0.9994"
Set-TextValue $ws.Range("E24") "  -0.16%  "

Set-TextValue $ws.Range("E25") "  -6.32%  "

Set-TextValue $ws.Range("D26") "9.267"
Set-TextValue $ws.Range("E26") "  -2.34%  "

Set-TextValue $ws.Range("D27") "164.84"
Set-TextValue $ws.Range("E27") "  +0.92%  "

Set-TextValue $ws.Range("E29") "  -0.84%  "

Set-TextValue $ws.Range("D30") "1.512"
Set-TextValue $ws.Range("E30") "  +5.38%  "

Set-TextValue $ws.Range("E31") "  +1.61%  "

Set-TextValue $ws.Range("D32") "1.532"
Set-TextValue $ws.Range("E32") "  -0.77%  "

Set-TextValue $ws.Range("D33") "4.272"
Set-TextValue $ws.Range("E33") "  +3.40%  "

Set-TextValue $ws.Range("D34") "0.05317"
Set-TextValue $ws.Range("E34") "  -2.89%  "

Set-TextValue $ws.Range("D35") "1.237"
Set-TextValue $ws.Range("E35") "  -0.68%  "

Set-TextValue $ws.Range("D36") "0.7543"
Set-TextValue $ws.Range("E36") "  -0.23%  "

Set-TextValue $ws.Range("D37") "1.001"
Set-TextValue $ws.Range("E37") "  -0.48%  "

Set-TextValue $ws.Range("D38") "2.694"
Set-TextValue $ws.Range("E38") "  +0.13%  "

Set-TextValue $ws.Range("D39") "0.01958"
Set-TextValue $ws.Range("E39") "  -0.56%  "

Set-TextValue $ws.Range("D40") "2.752"
Set-TextValue $ws.Range("E40") "  -1.30%  "

Set-TextValue $ws.Range("D41") "0.4514"
Set-TextValue $ws.Range("E41") "  +0.10%  "

Set-TextValue $ws.Range("D42") "1.113.15"
Set-TextValue $ws.Range("E42") "  +0.39%  "

Set-TextValue $ws.Range("D43") "6.060"
Set-TextValue $ws.Range("E43") "  -0.72%  "

Set-TextValue $ws.Range("D44") "72.37"
Set-TextValue $ws.Range("E44") "  -1.89%  "

Set-TextValue $ws.Range("D45") "0.8552"
Set-TextValue $ws.Range("E45") "  +0.24%  "

Set-TextValue $ws.Range("D46") "1.001"
Set-TextValue $ws.Range("E46") "  +0.05%  "

Set-TextValue $ws.Range("D47") "103.42"
Set-TextValue $ws.Range("E47") "  -0.26%  "

Set-TextValue $ws.Range("D48") "3.108"
Set-TextValue $ws.Range("E48") "  +3.10%  "

Set-TextValue $ws.Range("D49") "7.631"
Set-TextValue $ws.Range("E49") "  +0.27%  "

Set-TextValue $ws.Range("D50") "1.839"
Set-TextValue $ws.Range("E50") "  -2.17%  "

Set-TextValue $ws.Range("D51") "2.025.94"
Set-TextValue $ws.Range("E51") "  -7.67%  "
